$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @{
  "2" = @{
    "D" = $null
    "E" = $null
    "F" = $null
    "G" = $null
    "H" = $null
    "I" = $null
    "J" = $null
    "K" = $null
    "L" = $null
    "M" = $null
    "N" = $null
    "O" = $null
    "P" = $null
    "Q" = $null
    "R" = $null
    "S" = $null
    "T" = $null
    "U" = $null
    "V" = $null
    "W" = $null
    "X" = $null
    "Y" = $null
    "Z" = $null
    "AA" = $null
    "AB" = $null
    "AC" = $null
    "AD" = $null
    "AE" = $null
    "AF" = $null
    "AG" = 1300
    "AH" = 0.9399999999999999
    "AI" = $null
    "AJ" = 12000000
  }
  "3" = @{
    "D" = 8960
    "E" = 1046
    "F" = 1046
    "G" = 1462
    "H" = 1111
    "I" = 1111
    "J" = $null
    "K" = 16664
    "L" = 3095
    "M" = 13569
    "N" = 13569
    "O" = $null
    "P" = 600
    "Q" = 851
    "R" = -373
    "S" = -170
    "T" = 127
    "U" = 724
    "V" = $null
    "W" = 11.67
    "X" = 12.4
    "Y" = $null
    "Z" = $null
    "AA" = 22.81
    "AB" = 2175.42
    "AC" = 9259
    "AD" = 12.47
    "AE" = 114563
    "AF" = 1.01
    "AG" = 1500
    "AH" = 1.3
    "AI" = 15.99
    "AJ" = 12000000
  }
  "4" = @{
    "D" = 9694
    "E" = 1113
    "F" = 1113
    "G" = 1741
    "H" = 1295
    "I" = 1295
    "J" = $null
    "K" = 19030
    "L" = 4477
    "M" = 14553
    "N" = 14553
    "O" = $null
    "P" = 600
    "Q" = 928
    "R" = -1945
    "S" = 718
    "T" = 92
    "U" = 836
    "V" = 960
    "W" = 11.48
    "X" = 13.35
    "Y" = 9.210000000000001
    "Z" = 7.25
    "AA" = 30.76
    "AB" = 2360.69
    "AC" = 10788
    "AD" = 10.34
    "AE" = 123533
    "AF" = 0.9
    "AG" = 1500
    "AH" = 1.35
    "AI" = 13.67
    "AJ" = 12000000
  }
  "5" = @{
    "D" = 10431
    "E" = 1253
    "F" = 1253
    "G" = 1618
    "H" = 1225
    "I" = 1225
    "J" = $null
    "K" = 18754
    "L" = 3310
    "M" = 15444
    "N" = 15444
    "O" = $null
    "P" = 600
    "Q" = 279
    "R" = 1227
    "S" = -1258
    "T" = 101
    "U" = 179
    "V" = 0
    "W" = 12.02
    "X" = 11.74
    "Y" = 8.16
    "Z" = 6.48
    "AA" = 21.43
    "AB" = 2535.32
    "AC" = 10205
    "AD" = 11.81
    "AE" = 132142
    "AF" = 0.91
    "AG" = 1700
    "AH" = 1.41
    "AI" = 16.22
    "AJ" = 12000000
  }
  "6" = @{
    "D" = 10177
    "E" = 1123
    "F" = 1123
    "G" = 1856
    "H" = 1666
    "I" = 1666
    "K" = 25477
    "L" = 8632
    "M" = 16845
    "N" = 16845
    "P" = 600
    "Q" = 20
    "R" = 583
    "S" = -50
    "T" = 135
    "U" = -115
    "V" = 2306
    "W" = 11.04
    "X" = 16.37
    "Y" = 10.32
    "Z" = 7.53
    "AA" = 51.24
    "AB" = 2771.35
    "AC" = 13881
    "AD" = 7.13
    "AE" = 144125
    "AF" = 0.6899999999999999
    "AG" = $null
    "AH" = $null
    "AI" = 13.33
    "AJ" = 12000000
  }
  "7" = @{
    "D" = 21849
    "E" = 1292
    "G" = 1823
    "H" = 1407
    "I" = 1401
    "K" = 28341
    "L" = 10355
    "M" = 17986
    "N" = 17982
    "P" = 600
    "Q" = 1541
    "R" = -410
    "S" = 166
    "T" = 275
    "U" = 657
    "W" = 5.91
    "X" = 6.44
    "Y" = 8.050000000000001
    "Z" = 5.23
    "AA" = 57.57
    "AC" = 11678
    "AD" = 6.43
    "AE" = 153851
    "AF" = 0.49
    "AG" = 1988
    "AH" = 2.65
    "AI" = 17.02
  }
  "8" = @{
    "D" = 23108
    "E" = 1479
    "G" = 1988
    "H" = 1516
    "I" = 1521
    "K" = 30251
    "L" = 11002
    "M" = 19247
    "N" = 19233
    "P" = 600
    "Q" = 2209
    "R" = -547
    "S" = -397
    "T" = 297
    "U" = 971
    "W" = 6.4
    "X" = 6.56
    "Y" = 8.17
    "Z" = 5.17
    "AA" = 57.16
    "AC" = 12673
    "AD" = 5.93
    "AE" = 164560
    "AF" = 0.46
    "AG" = 2112
    "AH" = 2.81
    "AI" = 16.67
  }
  "9" = @{
    "D" = 24227
    "E" = 1639
    "G" = 2166
    "H" = 1645
    "I" = 1649
    "K" = 32877
    "L" = 12380
    "M" = 20497
    "N" = 20451
    "P" = 600
    "Q" = 2305
    "R" = -519
    "S" = -337
    "T" = 306
    "U" = 1129
    "W" = 6.76
    "X" = 6.79
    "Y" = 8.31
    "Z" = 5.21
    "AA" = 60.4
    "AC" = 13742
    "AD" = 5.47
    "AE" = 174981
    "AF" = 0.43
    "AG" = 2225
    "AH" = 2.96
    "AI" = 16.19
  }
}

foreach ($rowKey in $rowData.Keys) {
  $cols = $rowData[$rowKey]
  foreach ($colKey in $cols.Keys) {
    $addr = "${colKey}${rowKey}"
    $ws.Range($addr).Value2 = $cols[$colKey]
  }
}
